$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.155.42"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "1.837.03"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("D4").Value = "'0.9992"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'240.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.43%  "
$ws.Range("D6").Value = "'0.6857"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.07%  "
$ws.Range("D7").Value = "'0.9995"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.3012"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.65%  "
$ws.Range("D9").Value = "'0.07464"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.47%  "
$ws.Range("D10").Value = "'23.14"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.88%  "
$ws.Range("E11").Value = "  -2.04%  "
$ws.Range("D12").Value = "1.839.76"
$ws.Range("D13").Value = "'5.057"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.45%  "
$ws.Range("D14").Value = "'0.6817"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("D15").Value = "'87.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.14%  "
$ws.Range("D16").Value = "'6.153"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -7.32%  "
$ws.Range("D17").Value = "29.140.13"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").Value = "'0.000008173"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.87%  "
$ws.Range("D19").Value = "2.081.13"
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("D20").Value = "'227.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E21").Value = "  -1.96%  "
$ws.Range("D22").Value = "'0.9996"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "'7.394"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.76%  "
$ws.Range("D24").Value = "'0.9999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "'0.1454"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.23%  "
$ws.Range("D26").Value = "'160.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.67%  "
$ws.Range("D27").Value = "'8.757"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.92%  "
$ws.Range("E28").Value = "  -1.05%  "
$ws.Range("D29").Value = "'1.512"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.08%  "
$ws.Range("D30").Value = "'4.274"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.80%  "
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("D32").Value = "'1.200"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.26%  "
$ws.Range("D33").Value = "'0.05178"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.99%  "
$ws.Range("D34").Value = "'0.7658"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.37%  "
$ws.Range("D35").Value = "'1.843"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.54%  "
$ws.Range("D36").Value = "'1.134"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.47%  "
$ws.Range("D37").Value = "'2.674"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.61%  "
$ws.Range("D38").Value = "1.314.45"
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("D39").Value = "'0.01833"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.11%  "
$ws.Range("D40").Value = "'2.718"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.31%  "
$ws.Range("D41").Value = "'0.9341"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.48%  "
$ws.Range("D42").Value = "'5.792"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.75%  "
$ws.Range("D43").Value = "'104.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.88%  "
$ws.Range("D44").Value = "'0.9997"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("E45").Value = "  +0.37%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.982.33"
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'65.05"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.21%  "
$ws.Range("D48").Value = "'0.5198"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.41%  "
$ws.Range("D49").Value = "'9.543"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.01%  "
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("D51").Value = "'0.05936"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.89%  "
